$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1619000
$ws.Range("J86").Value = 3615.25
$ws.Range("L86").Value = 3615.25
$ws.Range("N86").Value = -5861.25
$ws.Range("H89").Value = 1619000
$ws.Range("J89").Value = 3615.25
$ws.Range("L89").Value = 18076.25
$ws.Range("N89").Value = -29308.25
$ws.Range("H96").Value = 2764.4546
$ws.Range("I96").Value = 465.7143
$ws.Range("K96").Value = 1397.1429
$ws.Range("M96").Value = -24.14289999999983
$ws.Range("H98").Value = 1184.6842
$ws.Range("I98").Value = 1184.6842
$ws.Range("K98").Value = 1184.6842
$ws.Range("M98").Value = 313.3158000000001
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 2831.1667
$ws.Range("J111").Value = 3149
$ws.Range("L111").Value = 9447
$ws.Range("N111").Value = -15581
$ws.Range("H112").Value = 3349.875
$ws.Range("I112").Value = 1575
$ws.Range("J112").Value = 5124.75
$ws.Range("K112").Value = 4725
$ws.Range("L112").Value = 15374.25
$ws.Range("M112").Value = -3617
$ws.Range("N112").Value = -17590.25
$ws.Range("H122").Value = 1184.6842
$ws.Range("I122").Value = 1184.6842
$ws.Range("K122").Value = 3554.0526
$ws.Range("M122").Value = -1104.0526
$ws.Range("H132").Value = 12474.4
$ws.Range("I132").Value = 4656.615
$ws.Range("K132").Value = 13969.845
$ws.Range("M132").Value = -11439.845
$ws.Range("H136").Value = 117440
$ws.Range("J136").Value = 117440
$ws.Range("L136").Value = 117440
$ws.Range("N136").Value = -127640
$ws.Range("H137").Value = 3836.303
$ws.Range("I137").Value = 1117.25
$ws.Range("K137").Value = 3351.75
$ws.Range("M137").Value = -801.75
$ws.Range("H138").Value = 3389.54
$ws.Range("J138").Value = 3930.2974
$ws.Range("L138").Value = 11790.8922
$ws.Range("N138").Value = -22070.8922
$ws.Range("H139").Value = 121155.8
$ws.Range("J139").Value = 121155.8
$ws.Range("L139").Value = 121155.8
$ws.Range("N139").Value = -131435.8
$ws.Range("H140").Value = 55875.152
$ws.Range("J140").Value = 54639
$ws.Range("L140").Value = 54639
$ws.Range("N140").Value = -64999

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3221.092
$ws.Range("I32").Value = 1802.4615
$ws.Range("J32").Value = 11603.909
$ws.Range("K32").Value = 1802.4615
$ws.Range("L32").Value = 11603.909
$ws.Range("M32").Value = -1515.4615
$ws.Range("N32").Value = -12177.909
$ws.Range("H74").Value = 22728866
$ws.Range("I74").Value = 41667970
$ws.Range("K74").Value = 41667970
$ws.Range("M74").Value = -41667096
$ws.Range("H77").Value = 22728866
$ws.Range("I77").Value = 41667970
$ws.Range("K77").Value = 208339850
$ws.Range("M77").Value = -208335482

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 5050
$ws.Range("J44").Value = 5050
$ws.Range("L44").Value = 5050
$ws.Range("N44").Value = -6044
$ws.Range("H107").Value = 1260.6765
$ws.Range("I107").Value = 1164.8572
$ws.Range("J107").Value = 1707.8334
$ws.Range("K107").Value = 1164.8572
$ws.Range("L107").Value = 1707.8334
$ws.Range("M107").Value = 755.1428000000001
$ws.Range("N107").Value = -5547.8334

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4416.1665
$ws.Range("I10").Value = 3832.3333
$ws.Range("K10").Value = 3832.3333
$ws.Range("M10").Value = -3693.3333
$ws.Range("H62").Value = 41463.75
$ws.Range("J62").Value = 63221.2
$ws.Range("L62").Value = 63221.2
$ws.Range("N62").Value = -64469.2
$ws.Range("H65").Value = 41463.75
$ws.Range("J65").Value = 63221.2
$ws.Range("L65").Value = 316106
$ws.Range("N65").Value = -322346
$ws.Range("H107").Value = 1818635.6
$ws.Range("I107").Value = 2597693.8
$ws.Range("K107").Value = 2597693.8
$ws.Range("M107").Value = -2595773.8
$ws.Range("H132").Value = 12823030
$ws.Range("I132").Value = 18520788
$ws.Range("J132").Value = 3075.125
$ws.Range("K132").Value = 55562364
$ws.Range("L132").Value = 9225.375
$ws.Range("M132").Value = -55559834
$ws.Range("N132").Value = -14285.375
$ws.Range("H134").Value = 2692.3777
$ws.Range("I134").Value = 2744.075
$ws.Range("J134").Value = 2278.8
$ws.Range("K134").Value = 8232.224999999999
$ws.Range("L134").Value = 6836.400000000001
$ws.Range("M134").Value = -5697.224999999999
$ws.Range("N134").Value = -11906.4
$ws.Range("H141").Value = 91822.5
$ws.Range("I141").Value = 39298.5
$ws.Range("J141").Value = 100576.5
$ws.Range("K141").Value = 39298.5
$ws.Range("L141").Value = 100576.5
$ws.Range("M141").Value = -34118.5
$ws.Range("N141").Value = -110936.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1460.25
$ws.Range("J129").Value = 1408
$ws.Range("L129").Value = 4224
$ws.Range("N129").Value = -14224
$ws.Range("H133").Value = 13637.071
$ws.Range("I133").Value = 4159.6665
$ws.Range("J133").Value = 20745.125
$ws.Range("K133").Value = 12478.9995
$ws.Range("L133").Value = 62235.375
$ws.Range("M133").Value = -7418.999500000002
$ws.Range("N133").Value = -72355.375
$ws.Range("H137").Value = 7348179.5
$ws.Range("J137").Value = 7872692.5
$ws.Range("L137").Value = 23618077.5
$ws.Range("N137").Value = -23628277.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2976778
$ws.Range("J107").Value = 745.25
$ws.Range("L107").Value = 745.25
$ws.Range("N107").Value = -4585.25
$ws.Range("H132").Value = 6521.7427
$ws.Range("I132").Value = 5743.231
$ws.Range("K132").Value = 17229.693
$ws.Range("M132").Value = -14699.693

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12502417
$ws.Range("I16").Value = 13335311
$ws.Range("J16").Value = 9002
$ws.Range("K16").Value = 13335311
$ws.Range("L16").Value = 9002
$ws.Range("M16").Value = -13335141
$ws.Range("N16").Value = -9342

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1198.5714
$ws.Range("I126").Value = 1198.3334
$ws.Range("K126").Value = 3595.0002
$ws.Range("M126").Value = -1125.0002
$ws.Range("H132").Value = 24161598
$ws.Range("I132").Value = 5559483
$ws.Range("J132").Value = 38470916
$ws.Range("K132").Value = 16678449
$ws.Range("L132").Value = 115412748
$ws.Range("M132").Value = -16675919
$ws.Range("N132").Value = -115417808
$ws.Range("H138").Value = 74333.336
$ws.Range("J138").Value = 74333.336
$ws.Range("L138").Value = 74333.336
$ws.Range("N138").Value = -84613.336
